$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.773.03'
$ws.Range("E2").Value = '  -2.08%  '
$ws.Range("D3").Value = '3.232.11'
$ws.Range("E3").Value = '  -1.35%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '''577.07'
$ws.Range("E5").Value = '  -1.22%  '
$ws.Range("D6").Value = '''172.61'
$ws.Range("E6").Value = '  -3.66%  '
$ws.Range("D7").Value = '''0.630'
$ws.Range("E7").Value = '  +0.67%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("D9").Value = '3.228.00'
$ws.Range("E9").Value = '  -1.44%  '
$ws.Range("E10").Value = '  -2.62%  '
$ws.Range("D11").Value = '''6.76'
$ws.Range("E11").Value = '  +0.43%  '
$ws.Range("E12").Value = '  -3.14%  '
$ws.Range("D13").Value = '3.794.43'
$ws.Range("E13").Value = '  -1.41%  '
$ws.Range("D15").Value = '64.855.30'
$ws.Range("E15").Value = '  -1.97%  '
$ws.Range("D16").Value = '''25.78'
$ws.Range("E16").Value = '  -2.11%  '
$ws.Range("D17").Value = '3.231.08'
$ws.Range("E17").Value = '  -1.08%  '
$ws.Range("E18").Value = '  -2.68%  '
$ws.Range("D19").Value = '''416.24'
$ws.Range("E19").Value = '  -4.19%  '
$ws.Range("D20").Value = '''5.38'
$ws.Range("E20").Value = '  -2.19%  '
$ws.Range("D21").Value = '''12.79'
$ws.Range("E21").Value = '  -2.92%  '
$ws.Range("D22").Value = '''7.21'
$ws.Range("E22").Value = '  -2.21%  '
$ws.Range("D23").Value = '''0.997'
$ws.Range("E23").Value = '  -0.22%  '
$ws.Range("D24").Value = '''70.56'
$ws.Range("E24").Value = '  -1.51%  '
$ws.Range("E25").Value = '  -0.61%  '
$ws.Range("E26").Value = '  +4.34%  '
$ws.Range("D27").Value = '''0.495'
$ws.Range("E27").Value = '  -1.89%  '
$ws.Range("D28").Value = '''0.0000111'
$ws.Range("E28").Value = '  -1.44%  '
$ws.Range("D29").Value = '''9.02'
$ws.Range("E29").Value = '  +2.11%  '
$ws.Range("E30").Value = '  -0.04%  '
$ws.Range("E31").Value = '  -4.33%  '
$ws.Range("D32").Value = '''21.81'
$ws.Range("E32").Value = '  -1.96%  '
$ws.Range("E33").Value = '  +0.08%  '
$ws.Range("D34").Value = '''4.98'
$ws.Range("E34").Value = '  -3.97%  '
$ws.Range("D35").Value = '''6.42'
$ws.Range("E35").Value = '  -2.76%  '
$ws.Range("E36").Value = '  -2.55%  '
$ws.Range("D37").Value = '''157.53'
$ws.Range("E37").Value = '  -0.27%  '
$ws.Range("E38").Value = '  -2.17%  '
$ws.Range("D39").Value = '2.826.49'
$ws.Range("E39").Value = '  +1.93%  '
$ws.Range("D40").Value = '''1.73'
$ws.Range("E40").Value = '  -2.87%  '
$ws.Range("D41").Value = '''25.42'
$ws.Range("E41").Value = '  -4.37%  '
$ws.Range("D42").Value = '''4.22'
$ws.Range("E42").Value = '  -2.49%  '
$ws.Range("D43").Value = '''39.49'
$ws.Range("E43").Value = '  -1.69%  '
$ws.Range("D44").Value = '''0.723'
$ws.Range("E44").Value = '  -6.51%  '
$ws.Range("E45").Value = '  -4.56%  '
$ws.Range("D46").Value = '''0.0630'
$ws.Range("E46").Value = '  -4.06%  '
$ws.Range("B47").Value = 'dogwifhat'
$ws.Range("C47").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D47").Value = '''2.17'
$ws.Range("E47").Value = '  -4.79%  '
$ws.Range("B48").Value = 'Bittensor'
$ws.Range("C48").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D48").Value = '''301.52'
$ws.Range("E48").Value = '  -5.70%  '
$ws.Range("D49").Value = '''21.97'
$ws.Range("E49").Value = '  -5.30%  '
$ws.Range("D50").Value = '''0.0263'
$ws.Range("E50").Value = '  -1.32%  '
$ws.Range("E51").Value = '  -1.08%  '